$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F holds "date" values (serial date numbers). Shift each of the six
# dates in F2:F7 forward by 2 days, keeping the existing date number format.
$ws.Range("F2").Value = 44636
$ws.Range("F3").Value = 44635
$ws.Range("F4").Value = 44634
$ws.Range("F5").Value = 44633
$ws.Range("F6").Value = 44632
$ws.Range("F7").Value = 44631
